$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 249.66667
$ws.Range("J12").Value = 249
$ws.Range("L12").Value = 249
$ws.Range("N12").Value = -589

$ws.Range("H40").Value = 1301.7142
$ws.Range("I40").Value = 1028
$ws.Range("J40").Value = 1666.6666
$ws.Range("K40").Value = 1028
$ws.Range("L40").Value = 1666.6666
$ws.Range("M40").Value = -853
$ws.Range("N40").Value = -2016.6666

$ws.Range("H116").Value = 3147.5
$ws.Range("J116").Value = 3177
$ws.Range("L116").Value = 3177
$ws.Range("N116").Value = -10061

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2353.037
$ws.Range("I2").Value = 1395.4
$ws.Range("J2").Value = 3550.0833
$ws.Range("K2").Value = 1395.4
$ws.Range("L2").Value = 3550.0833
$ws.Range("M2").Value = -1282.4
$ws.Range("N2").Value = -3776.0833

$ws.Range("H32").Value = 3938.3333
$ws.Range("I32").Value = 2798.5454
$ws.Range("K32").Value = 2798.5454
$ws.Range("M32").Value = -2511.5454

$ws.Range("H74").Value = 1052.5294
$ws.Range("I74").Value = 1052.5294
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1052.5294
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -178.5293999999999
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1052.5294
$ws.Range("I77").Value = 1052.5294
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5262.646999999999
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -894.646999999999
$ws.Range("N77").ClearContents()

$ws.Range("H116").Value = 2353.037
$ws.Range("I116").Value = 1395.4
$ws.Range("J116").Value = 3550.0833
$ws.Range("K116").Value = 1395.4
$ws.Range("L116").Value = 3550.0833
$ws.Range("M116").Value = 898.5999999999999
$ws.Range("N116").Value = -8138.0833

$ws.Range("H122").Value = 3128.3333
$ws.Range("I122").Value = 2335.6
$ws.Range("J122").Value = 4119.25
$ws.Range("K122").Value = 7006.799999999999
$ws.Range("L122").Value = 12357.75
$ws.Range("M122").Value = -4556.799999999999
$ws.Range("N122").Value = -17257.75

$ws.Range("H132").Value = 1860.4
$ws.Range("I132").Value = 1030.3334
$ws.Range("J132").Value = 3105.5
$ws.Range("K132").Value = 3091.0002
$ws.Range("L132").Value = 9316.5
$ws.Range("M132").Value = -561.0001999999999
$ws.Range("N132").Value = -14376.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2353.037
$ws.Range("I3").Value = 1395.4
$ws.Range("J3").Value = 3550.0833
$ws.Range("K3").Value = 1395.4
$ws.Range("L3").Value = 3550.0833
$ws.Range("M3").Value = -1281.4
$ws.Range("N3").Value = -3778.0833

$ws.Range("H81").Value = 17250
$ws.Range("J81").Value = 17250
$ws.Range("L81").Value = 17250
$ws.Range("N81").Value = -19372

$ws.Range("H84").Value = 17250
$ws.Range("J84").Value = 17250
$ws.Range("L84").Value = 51750
$ws.Range("N84").Value = -62358

$ws.Range("H134").Value = 2015.875
$ws.Range("I134").Value = 2082.6667
$ws.Range("J134").Value = 1014
$ws.Range("K134").Value = 6248.000100000001
$ws.Range("L134").Value = 3042
$ws.Range("M134").Value = -3713.000100000001
$ws.Range("N134").Value = -8112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1539.3334
$ws.Range("I31").Value = 1506.9231
$ws.Range("J31").Value = 1750
$ws.Range("K31").Value = 1506.9231
$ws.Range("L31").Value = 1750
$ws.Range("M31").Value = -1211.9231
$ws.Range("N31").Value = -2340

$ws.Range("H34").Value = 1539.3334
$ws.Range("I34").Value = 1506.9231
$ws.Range("J34").Value = 1750
$ws.Range("K34").Value = 1506.9231
$ws.Range("L34").Value = 1750
$ws.Range("M34").Value = -1304.9231
$ws.Range("N34").Value = -2154

$ws.Range("H132").Value = 2227.862
$ws.Range("I132").Value = 2267.074
$ws.Range("K132").Value = 6801.222
$ws.Range("M132").Value = -4271.222

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 502.33334
$ws.Range("I4").Value = 466.0909
$ws.Range("K4").Value = 1398.2727
$ws.Range("M4").Value = -1286.2727

$ws.Range("H81").Value = 9000
$ws.Range("I81").Value = 9000
$ws.Range("K81").Value = 27000
$ws.Range("M81").Value = -25877

$ws.Range("H84").Value = 9000
$ws.Range("I84").Value = 9000
$ws.Range("K84").Value = 81000
$ws.Range("M84").Value = -75384

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 25000
$ws.Range("J15").Value = 25000
$ws.Range("L15").Value = 25000
$ws.Range("N15").Value = -25576

$ws.Range("H81").Value = 25000
$ws.Range("J81").Value = 25000
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -26996

$ws.Range("H84").Value = 25000
$ws.Range("J84").Value = 25000
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -84984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5120
$ws.Range("I46").Value = 4900
$ws.Range("K46").Value = 4900
$ws.Range("M46").Value = -4712

$ws.Range("H132").Value = 1989.6666
$ws.Range("I132").Value = 1858.25
$ws.Range("J132").Value = 2252.5
$ws.Range("K132").Value = 5574.75
$ws.Range("L132").Value = 6757.5
$ws.Range("M132").Value = -3044.75
$ws.Range("N132").Value = -11817.5

$ws.Range("H136").Value = 3099.75
$ws.Range("I136").Value = 3099.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9299.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6749.25
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 28578.5
$ws.Range("I76").Value = 7157
$ws.Range("J76").Value = 50000
$ws.Range("K76").Value = 7157
$ws.Range("L76").Value = 50000
$ws.Range("M76").Value = -6842
$ws.Range("N76").Value = -50630

$ws.Range("H79").Value = 28578.5
$ws.Range("I79").Value = 7157
$ws.Range("J79").Value = 50000
$ws.Range("K79").Value = 7157
$ws.Range("L79").Value = 50000
$ws.Range("M79").Value = -6065
$ws.Range("N79").Value = -52184

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
